$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2523.875
$ws.Range("I51").Value = 2067
$ws.Range("J51").Value = 2798
$ws.Range("K51").Value = 2067
$ws.Range("L51").Value = 2798
$ws.Range("M51").Value = -1583
$ws.Range("N51").Value = -3766

$ws.Range("H69").Value = 18520948
$ws.Range("I69").Value = 5000
$ws.Range("J69").Value = 18870306
$ws.Range("K69").Value = 15000
$ws.Range("L69").Value = 56610918
$ws.Range("M69").Value = -14126
$ws.Range("N69").Value = -56612666

$ws.Range("H72").Value = 18520948
$ws.Range("I72").Value = 5000
$ws.Range("J72").Value = 18870306
$ws.Range("K72").Value = 45000
$ws.Range("L72").Value = 169832754
$ws.Range("M72").Value = -40632
$ws.Range("N72").Value = -169841490

$ws.Range("H74").Value = 3596.2964
$ws.Range("I74").Value = 4333.3335
$ws.Range("J74").Value = 3504.1667
$ws.Range("K74").Value = 4333.3335
$ws.Range("L74").Value = 3504.1667
$ws.Range("M74").Value = -3397.3335
$ws.Range("N74").Value = -5376.1667

$ws.Range("H77").Value = 3596.2964
$ws.Range("I77").Value = 4333.3335
$ws.Range("J77").Value = 3504.1667
$ws.Range("K77").Value = 21666.6675
$ws.Range("L77").Value = 17520.8335
$ws.Range("M77").Value = -16986.6675
$ws.Range("N77").Value = -26880.8335

$ws.Range("H81").Value = 70000
$ws.Range("J81").Value = 70000
$ws.Range("L81").Value = 70000
$ws.Range("N81").Value = -71996

$ws.Range("H84").Value = 70000
$ws.Range("J84").Value = 70000
$ws.Range("L84").Value = 210000
$ws.Range("N84").Value = -219984

$ws.Range("H86").Value = 1746.7894
$ws.Range("I86").Value = 1745.4615
$ws.Range("J86").Value = 1749.6666
$ws.Range("K86").Value = 1745.4615
$ws.Range("L86").Value = 1749.6666
$ws.Range("M86").Value = -622.4614999999999
$ws.Range("N86").Value = -3995.6666

$ws.Range("H89").Value = 1746.7894
$ws.Range("I89").Value = 1745.4615
$ws.Range("J89").Value = 1749.6666
$ws.Range("K89").Value = 8727.307499999999
$ws.Range("L89").Value = 8748.333000000001
$ws.Range("M89").Value = -3111.307499999999
$ws.Range("N89").Value = -19980.333

$ws.Range("H106").Value = 1617.5834
$ws.Range("I106").Value = 1380.5
$ws.Range("J106").Value = 2803
$ws.Range("K106").Value = 1380.5
$ws.Range("L106").Value = 2803
$ws.Range("M106").Value = -749.5
$ws.Range("N106").Value = -4065

$ws.Range("H111").Value = 55556000
$ws.Range("I111").Value = 62500460
$ws.Range("J111").Value = 300
$ws.Range("K111").Value = 187501380
$ws.Range("L111").Value = 900
$ws.Range("M111").Value = -187498313
$ws.Range("N111").Value = -7034

$ws.Range("H129").Value = 703.5

$ws.Range("H137").Value = 10389539
$ws.Range("I137").Value = 16968216
$ws.Range("J137").Value = 2155.7368
$ws.Range("K137").Value = 50904648
$ws.Range("L137").Value = 6467.2104
$ws.Range("M137").Value = -50902098
$ws.Range("N137").Value = -11567.2104

$ws.Range("H138").Value = 2992.2
$ws.Range("I138").Value = 1183.875
$ws.Range("J138").Value = 3336.6428
$ws.Range("K138").Value = 3551.625
$ws.Range("L138").Value = 10009.9284
$ws.Range("M138").Value = 1588.375
$ws.Range("N138").Value = -20289.9284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10768182
$ws.Range("I32").Value = 14501405
$ws.Range("J32").Value = 35166.5
$ws.Range("K32").Value = 14501405
$ws.Range("L32").Value = 35166.5
$ws.Range("M32").Value = -14501118
$ws.Range("N32").Value = -35740.5

$ws.Range("H97").Value = 779.9231
$ws.Range("I97").Value = 567.8
$ws.Range("J97").Value = 1487
$ws.Range("K97").Value = 567.8
$ws.Range("L97").Value = 1487
$ws.Range("M97").Value = -71.79999999999995
$ws.Range("N97").Value = -2479

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2017.5
$ws.Range("I105").Value = 1799.875
$ws.Range("J105").Value = 2452.75
$ws.Range("K105").Value = 1799.875
$ws.Range("L105").Value = 2452.75
$ws.Range("M105").Value = -52.875
$ws.Range("N105").Value = -5946.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 516305.7
$ws.Range("I31").Value = 1943.7037
$ws.Range("J31").Value = 1071816.6
$ws.Range("K31").Value = 1943.7037
$ws.Range("L31").Value = 1071816.6
$ws.Range("M31").Value = -1648.7037
$ws.Range("N31").Value = -1072406.6

$ws.Range("H34").Value = 516305.7
$ws.Range("I34").Value = 1943.7037
$ws.Range("J34").Value = 1071816.6
$ws.Range("K34").Value = 1943.7037
$ws.Range("L34").Value = 1071816.6
$ws.Range("M34").Value = -1741.7037
$ws.Range("N34").Value = -1072220.6

$ws.Range("H118").Value = 34750
$ws.Range("J118").Value = 34750
$ws.Range("L118").Value = 34750
$ws.Range("N118").Value = -38064

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 845.2692
$ws.Range("I131").Value = 327.77777
$ws.Range("J131").Value = 953.5814
$ws.Range("K131").Value = 983.33331
$ws.Range("L131").Value = 2860.7442
$ws.Range("M131").Value = 4056.66669
$ws.Range("N131").Value = -12940.7442

$ws.Range("H137").Value = 3390.25
$ws.Range("I137").Value = 2607.6924
$ws.Range("J137").Value = 4315.091
$ws.Range("K137").Value = 7823.0772
$ws.Range("L137").Value = 12945.273
$ws.Range("M137").Value = -2723.0772
$ws.Range("N137").Value = -23145.273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2961.3333
$ws.Range("I80").Value = 2536.3635
$ws.Range("J80").Value = 3320.923
$ws.Range("K80").Value = 2536.3635
$ws.Range("L80").Value = 3320.923
$ws.Range("M80").Value = -1538.3635
$ws.Range("N80").Value = -5316.923

$ws.Range("H83").Value = 2961.3333
$ws.Range("I83").Value = 2536.3635
$ws.Range("J83").Value = 3320.923
$ws.Range("K83").Value = 12681.8175
$ws.Range("L83").Value = 16604.615
$ws.Range("M83").Value = -7689.817499999999
$ws.Range("N83").Value = -26588.615

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1766.4584
$ws.Range("I82").Value = 1660.4
$ws.Range("J82").Value = 1842.2142
$ws.Range("K82").Value = 1660.4
$ws.Range("L82").Value = 1842.2142
$ws.Range("M82").Value = -1299.4
$ws.Range("N82").Value = -2564.2142

$ws.Range("H85").Value = 1766.4584
$ws.Range("I85").Value = 1660.4
$ws.Range("J85").Value = 1842.2142
$ws.Range("K85").Value = 1660.4
$ws.Range("L85").Value = 1842.2142
$ws.Range("M85").Value = -412.4000000000001
$ws.Range("N85").Value = -4338.2142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4277.778
$ws.Range("I81").Value = 2916.6667
$ws.Range("J81").Value = 7000
$ws.Range("K81").Value = 5833.3334
$ws.Range("L81").Value = 14000
$ws.Range("M81").Value = -4772.3334
$ws.Range("N81").Value = -16122

$ws.Range("H84").Value = 4277.778
$ws.Range("I84").Value = 2916.6667
$ws.Range("J84").Value = 7000
$ws.Range("K84").Value = 29166.667
$ws.Range("L84").Value = 70000
$ws.Range("M84").Value = -23862.667
$ws.Range("N84").Value = -80608
